# Apply updated Betfair back/lay odds values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q4").Value = 1.67

$ws.Range("AE6").Value = 85
$ws.Range("AN6").Value = 28

$ws.Range("G8").Value = 2.76
$ws.Range("J8").Value = 3.6

$ws.Range("G9").Value = 2.24
$ws.Range("I9").Value = 4.7
$ws.Range("L9").Value = 1.43
$ws.Range("N9").Value = 3.2
$ws.Range("W9").Value = 1.8

$ws.Range("H10").Value = 2.72
$ws.Range("I10").Value = 2.94
$ws.Range("K10").Value = 3.3

$ws.Range("U11").Value = 1.68
$ws.Range("X11").Value = 23
$ws.Range("Z11").Value = 170
$ws.Range("AC11").Value = 14
$ws.Range("AH11").Value = 44
$ws.Range("AJ11").Value = 9

$ws.Range("F12").Value = 2.82
$ws.Range("G12").Value = 2.88
$ws.Range("H12").Value = 2.82
$ws.Range("I12").Value = 2.88
$ws.Range("J12").Value = 3.3
$ws.Range("K12").Value = 3.35
$ws.Range("N12").Value = 3.45
